$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Förändrad" (changed) date column C for rows 2-15: 45184 -> 45186 ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# --- Add friendly display text as 2nd HYPERLINK() argument for rows 2 and 3 ---
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/artfynd/A 33062-2023.xlsx", "A 33062-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/kartor/A 33062-2023.png", "A 33062-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/klagomål/A 33062-2023.docx", "A 33062-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/klagomålsmail/A 33062-2023.docx", "A 33062-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/tillsyn/A 33062-2023.docx", "A 33062-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/tillsynsmail/A 33062-2023.docx", "A 33062-2023")'

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/artfynd/A 1061-2022.xlsx", "A 1061-2022")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/kartor/A 1061-2022.png", "A 1061-2022")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/klagomål/A 1061-2022.docx", "A 1061-2022")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/klagomålsmail/A 1061-2022.docx", "A 1061-2022")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/tillsyn/A 1061-2022.docx", "A 1061-2022")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_UPPLANDS_VASBY/tillsynsmail/A 1061-2022.docx", "A 1061-2022")'
